$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G
$data = @{
    2  = @(-0.06304482099587455, 1.284473961807544, 3.268916597491724, 1.808014545708005, 1.824894733127932, 51)
    3  = @(0.3848683813723368, 1.606026431035804, 7.828542112970349, 2.797953200639773, 2.799493110457711, 50)
    4  = @(0.8746662378248481, 1.947065274388025, 10.77140143740715, 3.281981328010132, 3.196064181409324, 49)
    5  = @(0.5898522373503049, 1.91023936529586, 11.52381436931002, 3.394674412857589, 3.378412877893673, 48)
    6  = @(0.6796283455196681, 2.053046112058658, 12.25639243636365, 3.500913086091062, 3.471440604892087, 47)
    7  = @(0.4620040488475629, 1.878470610035362, 11.8468855781181, 3.441930501639757, 3.456566949417224, 38)
    8  = @(0.5682316542479313, 1.945817767357757, 12.30934837299113, 3.508468094908536, 3.50990290671253, 37)
    9  = @(0.5226750354899881, 2.444751333587228, 20.24450866695628, 4.499389810513897, 4.585023615187054, 20)
    10 = @(-0.5732758996450275, 2.122935804733599, 8.802873317106101, 2.96696365281176, 3.02991975663685, 13)
    11 = @(-0.3393998147840875, 2.457134838492891, 7.456238556435196, 2.730611388761718, 3.029242133389085, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
